$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 27.5514752002093
$ws.Range("E2").Value = 27.25054931640625
$ws.Range("F2").Value = 29.03474467798217
$ws.Range("G2").Value = 26.26882171544944
$ws.Range("H2").Value = 5817000000
$ws.Range("I2").Value = "GOOGL"

$ws.Range("D3").Value = 26.99680391970704
$ws.Range("E3").Value = 32.64981460571289
$ws.Range("F3").Value = 35.42219363114338
$ws.Range("G3").Value = 26.80313992323004
$ws.Range("H3").Value = 5817000000
$ws.Range("I3").Value = "GOOGL"

$ws.Range("D4").Value = 31.67850840614184
$ws.Range("E4").Value = 36.61694717407227
$ws.Range("F4").Value = 37.36727286981908
$ws.Range("G4").Value = 31.33984461281439
$ws.Range("H4").Value = 5817000000
$ws.Range("I4").Value = "GOOGL"

$ws.Range("D5").Value = 37.84895204770889
$ws.Range("E5").Value = 37.80674362182617
$ws.Range("F5").Value = 38.19655307766485
$ws.Range("G5").Value = 34.15343974843713
$ws.Range("H5").Value = 5817000000
$ws.Range("I5").Value = "GOOGL"

$ws.Range("D6").Value = 37.59867654451598
$ws.Range("E6").Value = 35.15155410766602
$ws.Range("F6").Value = 39.27660247423416
$ws.Range("G6").Value = 34.91915584266965
$ws.Range("H6").Value = 5817000000
$ws.Range("I6").Value = "GOOGL"

$ws.Range("D7").Value = 35.01350745892713
$ws.Range("E7").Value = 39.29597091674805
$ws.Range("F7").Value = 39.9216527452363
$ws.Range("G7").Value = 34.71059637860652
$ws.Range("H7").Value = 5817000000
$ws.Range("I7").Value = "GOOGL"

$ws.Range("D8").Value = 39.85263024647205
$ws.Range("E8").Value = 40.21761322021485
$ws.Range("F8").Value = 41.66264846949318
$ws.Range("G8").Value = 39.53879671665068
$ws.Range("H8").Value = 5817000000
$ws.Range("I8").Value = "GOOGL"

$ws.Range("D9").Value = 39.75679408009914
$ws.Range("E9").Value = 40.72859191894531
$ws.Range("F9").Value = 43.05305985144016
$ws.Range("G9").Value = 39.57157480996585
$ws.Range("H9").Value = 5817000000
$ws.Range("I9").Value = "GOOGL"

$ws.Range("D10").Value = 42.14681153922361
$ws.Range("E10").Value = 45.90936279296875
$ws.Range("F10").Value = 46.47446170968723
$ws.Range("G10").Value = 41.44415731931716
$ws.Range("H10").Value = 5817000000
$ws.Range("I10").Value = "GOOGL"

$ws.Range("D11").Value = 46.34137550910245
$ws.Range("E11").Value = 46.951171875
$ws.Range("F11").Value = 49.9648859350205
$ws.Range("G11").Value = 45.45200942440623
$ws.Range("H11").Value = 5817000000
$ws.Range("I11").Value = "GOOGL"

$ws.Range("D12").Value = 48.44835115601325
$ws.Range("E12").Value = 51.29819488525391
$ws.Range("F12").Value = 52.81672088405051
$ws.Range("G12").Value = 47.76804190074551
$ws.Range("H12").Value = 5817000000
$ws.Range("I12").Value = "GOOGL"

$ws.Range("D13").Value = 52.29034440233525
$ws.Range("E13").Value = 58.70609283447266
$ws.Range("F13").Value = 59.48968968097881
$ws.Range("G13").Value = 52.29034440233525
$ws.Range("H13").Value = 5817000000
$ws.Range("I13").Value = "GOOGL"

$ws.Range("D14").Value = 51.0290464927241
$ws.Range("E14").Value = 50.58014297485352
$ws.Range("F14").Value = 54.49960877603787
$ws.Range("G14").Value = 49.37197652866953
$ws.Range("H14").Value = 5817000000
$ws.Range("I14").Value = "GOOGL"

$ws.Range("D15").Value = 55.3854983032786
$ws.Range("E15").Value = 60.94068145751953
$ws.Range("F15").Value = 64.12968487494881
$ws.Range("G15").Value = 54.95099514799846
$ws.Range("H15").Value = 5817000000
$ws.Range("I15").Value = "GOOGL"

$ws.Range("D16").Value = 60.2345616698925
$ws.Range("E16").Value = 54.15548706054688
$ws.Range("F16").Value = 60.80661640223318
$ws.Range("G16").Value = 50.01504511365445
$ws.Range("H16").Value = 5817000000
$ws.Range("I16").Value = "GOOGL"

$ws.Range("D17").Value = 51.00819248996493
$ws.Range("E17").Value = 55.90888977050781
$ws.Range("F17").Value = 55.99728083126573
$ws.Range("G17").Value = 50.76834587925548
$ws.Range("H17").Value = 5817000000
$ws.Range("I17").Value = "GOOGL"

$ws.Range("D18").Value = 58.97027330396304
$ws.Range("E18").Value = 59.5373649597168
$ws.Range("F18").Value = 64.40429834230383
$ws.Range("G18").Value = 58.74483083432825
$ws.Range("H18").Value = 5817000000
$ws.Range("I18").Value = "GOOGL"

$ws.Range("D19").Value = 54.67489973756347
$ws.Range("E19").Value = 60.49277496337891
$ws.Range("F19").Value = 62.9850838491221
$ws.Range("G19").Value = 54.39135579650275
$ws.Range("H19").Value = 5817000000
$ws.Range("I19").Value = "GOOGL"

$ws.Range("D20").Value = 60.7058067077664
$ws.Range("E20").Value = 62.50886917114258
$ws.Range("F20").Value = 64.51701814098398
$ws.Range("G20").Value = 57.75863484554971
$ws.Range("H20").Value = 5817000000
$ws.Range("I20").Value = "GOOGL"

$ws.Range("D21").Value = 66.95868268115149
$ws.Range("E21").Value = 71.14828491210938
$ws.Range("F21").Value = 74.51506337959323
$ws.Range("G21").Value = 66.86333959811857
$ws.Range("H21").Value = 5817000000
$ws.Range("I21").Value = "GOOGL"

$ws.Range("D22").Value = 55.8150324836901
$ws.Range("E22").Value = 66.87375640869141
$ws.Range("F22").Value = 67.54165041051544
$ws.Range("G22").Value = 53.38578835720912
$ws.Range("H22").Value = 5817000000
$ws.Range("I22").Value = "GOOGL"

$ws.Range("D23").Value = 70.47244007306874
$ws.Range("E23").Value = 73.88787841796875
$ws.Range("F23").Value = 78.80893863315268
$ws.Range("G23").Value = 70.22464507933984
$ws.Range("H23").Value = 5817000000
$ws.Range("I23").Value = "GOOGL"

$ws.Range("D24").Value = 73.67486378450468
$ws.Range("E24").Value = 80.25199890136719
$ws.Range("F24").Value = 83.49017093425739
$ws.Range("G24").Value = 71.17063578867412
$ws.Range("H24").Value = 5817000000
$ws.Range("I24").Value = "GOOGL"

$ws.Range("D25").Value = 87.39721568170974
$ws.Range("E25").Value = 90.7421417236328
$ws.Range("F25").Value = 95.94227611731162
$ws.Range("G25").Value = 84.22410116544214
$ws.Range("H25").Value = 5817000000
$ws.Range("I25").Value = "GOOGL"

$ws.Range("D26").Value = 103.8959189449111
$ws.Range("E26").Value = 116.8689422607422
$ws.Range("F26").Value = 120.7362660018762
$ws.Range("G26").Value = 103.8551994266863
$ws.Range("H26").Value = 5817000000
$ws.Range("I26").Value = "GOOGL"

$ws.Range("D27").Value = 120.8911910109136
$ws.Range("E27").Value = 133.8036346435547
$ws.Range("F27").Value = 137.3496726115578
$ws.Range("G27").Value = 120.699019745284
$ws.Range("H27").Value = 5817000000
$ws.Range("I27").Value = "GOOGL"

$ws.Range("D28").Value = 133.5265437133086
$ws.Range("E28").Value = 147.0319061279297
$ws.Range("F28").Value = 147.6317575184919
$ws.Range("G28").Value = 130.1523247071745
$ws.Range("H28").Value = 5817000000
$ws.Range("I28").Value = "GOOGL"

$ws.Range("D29").Value = 144.0613768134646
$ws.Range("E29").Value = 134.3766784667969
$ws.Range("F29").Value = 145.4815893005235
$ws.Range("G29").Value = 123.6471843017847
$ws.Range("H29").Value = 5817000000
$ws.Range("I29").Value = "GOOGL"

$ws.Range("D30").Value = 138.5444477205988
$ws.Range("E30").Value = 113.3278656005859
$ws.Range("F30").Value = 142.7276017833853
$ws.Range("G30").Value = 111.9642769129902
$ws.Range("H30").Value = 5817000000
$ws.Range("I30").Value = "GOOGL"

$ws.Range("D31").Value = 107.1936726298365
$ws.Range("E31").Value = 115.5232238769531
$ws.Range("F31").Value = 118.8651719734249
$ws.Range("G31").Value = 103.3571346729859
$ws.Range("H31").Value = 5817000000
$ws.Range("I31").Value = "GOOGL"

$ws.Range("D32").Value = 96.09720520386958
$ws.Range("E32").Value = 93.86261749267578
$ws.Range("F32").Value = 104.1019925357299
$ws.Range("G32").Value = 91.17118164756064
$ws.Range("H32").Value = 5817000000
$ws.Range("I32").Value = "GOOGL"

$ws.Range("D33").Value = 88.97631052391687
$ws.Range("E33").Value = 98.16294860839844
$ws.Range("F33").Value = 99.6328140358543
$ws.Range("G33").Value = 84.27871496823627
$ws.Range("H33").Value = 5817000000
$ws.Range("I33").Value = "GOOGL"

$ws.Range("D34").Value = 101.6886459843981
$ws.Range("E34").Value = 106.604736328125
$ws.Range("F34").Value = 108.4222029699671
$ws.Range("G34").Value = 101.2317978123801
$ws.Range("H34").Value = 5817000000
$ws.Range("I34").Value = "GOOGL"

$ws.Range("D35").Value = 118.423215877887
$ws.Range("E35").Value = 131.8108825683594
$ws.Range("F35").Value = 132.8238999142001
$ws.Range("G35").Value = 114.5598625925951
$ws.Range("H35").Value = 5817000000
$ws.Range("I35").Value = "GOOGL"

$ws.Range("D36").Value = 130.3112374168464
$ws.Range("E36").Value = 123.2300720214844
$ws.Range("F36").Value = 140.2526649298119
$ws.Range("G36").Value = 119.3865782259734
$ws.Range("H36").Value = 5817000000
$ws.Range("I36").Value = "GOOGL"

$ws.Range("D37").Value = 137.6009493970292
$ws.Range("E37").Value = 139.1403350830078
$ws.Range("F37").Value = 152.7266211780627
$ws.Range("G37").Value = 134.2242299641942
$ws.Range("H37").Value = 5817000000
$ws.Range("I37").Value = "GOOGL"

$ws.Range("D38").Value = 149.6577823870935
$ws.Range("E38").Value = 161.6649627685547
$ws.Range("F38").Value = 173.5132506603971
$ws.Range("G38").Value = 148.5752524773724
$ws.Range("H38").Value = 5817000000
$ws.Range("I38").Value = "GOOGL"

$ws.Range("D39").Value = 181.9848987390781
$ws.Range("E39").Value = 170.5605010986328
$ws.Range("F39").Value = 190.6551088124981
$ws.Range("G39").Value = 163.1232160954986
$ws.Range("H39").Value = 5817000000
$ws.Range("I39").Value = "GOOGL"

$ws.Range("D40").Value = 166.9537217165885
$ws.Range("E40").Value = 170.3587036132812
$ws.Range("F40").Value = 181.2208044470338
$ws.Range("G40").Value = 159.0386310205387
$ws.Range("H40").Value = 5817000000
$ws.Range("I40").Value = "GOOGL"

$ws.Range("D41").Value = 190.0304616757511
$ws.Range("E41").Value = 203.3570251464844
$ws.Range("F41").Value = 204.8122722439576
$ws.Range("G41").Value = 186.7511594827838
$ws.Range("H41").Value = 5817000000
$ws.Range("I41").Value = "GOOGL"

$ws.Range("D42").Value = 153.2971413880324
$ws.Range("E42").Value = 158.4662628173828
$ws.Range("F42").Value = 165.7509239001066
$ws.Range("G42").Value = 140.2346554932263
$ws.Range("H42").Value = 5817000000
$ws.Range("I42").Value = "GOOGL"

$ws.Range("D43").Value = 175.5829674872698
$ws.Range("E43").Value = 191.728515625
$ws.Range("F43").Value = 197.7731125063629
$ws.Range("G43").Value = 172.6156202045127
$ws.Range("H43").Value = 5817000000
$ws.Range("I43").Value = "GOOGL"
